# Apply the "Add files via upload" edits to ImageConversionTestCases.xlsx
# (Sheet1): update a handful of result/description cells, resize row 22,
# and refresh the current selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Content edits -------------------------------------------------------

# Test 4 (calcHori invalid width) & Test 7 (calcHori invalid height):
# expected result changed from "Return -1" to "Return 0".
$ws.Range("F11").Value = "Return 0"
$ws.Range("F14").Value = "Return 0"

# Test 10: test-name tweak (trailing period + space added).
$ws.Range("C18").Value = "Load image with valid file format and non-empty. "

# Test 13: test name expanded with the concrete test parameters.
$ws.Range("C22").Value = "Load image with valid file format and non-empty. For test, full black, w = 100 and h = 100"

# Row 22 grew taller (45 -> 60) to fit the longer wrapped text above.
$ws.Rows.Item(22).RowHeight = 60

# --- View state ------------------------------------------------------------
# Selection moved from D28 to the D22:D23 block, and the window scrolled
# back up so row 19 is the first visible row.
$ws.Range("D22:D23").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
